# Auto-generated COM-interop script to apply Famfrit_Profits.xlsx price-refresh diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Cells.Item(32, 8).Value = 0  # H32: 1581.1333 -> 0
$ws.Cells.Item(32, 9).Value = 0  # I32: 948.75 -> 0
$ws.Cells.Item(32, 10).Value = 0  # J32: 1811.091 -> 0
$ws.Cells.Item(32, 11).Value = 0  # K32: 948.75 -> 0
$ws.Cells.Item(32, 12).Value = 0  # L32: 1811.091 -> 0
$ws.Cells.Item(32, 13).ClearContents()  # M32: -622.75 -> (removed)
$ws.Cells.Item(32, 14).ClearContents()  # N32: -2463.091 -> (removed)
# Row 74
$ws.Cells.Item(74, 8).Value = 4563.364  # H74: 4198.9414 -> 4563.364
$ws.Cells.Item(74, 9).Value = 3366.1667  # I74: 3448.5 -> 3366.1667
$ws.Cells.Item(74, 11).Value = 3366.1667  # K74: 3448.5 -> 3366.1667
$ws.Cells.Item(74, 13).Value = -2430.1667  # M74: -2512.5 -> -2430.1667
# Row 77
$ws.Cells.Item(77, 8).Value = 4563.364  # H77: 4198.9414 -> 4563.364
$ws.Cells.Item(77, 9).Value = 3366.1667  # I77: 3448.5 -> 3366.1667
$ws.Cells.Item(77, 11).Value = 16830.8335  # K77: 17242.5 -> 16830.8335
$ws.Cells.Item(77, 13).Value = -12150.8335  # M77: -12562.5 -> -12150.8335
# Row 116
$ws.Cells.Item(116, 8).Value = 3101.4736  # H116: 3054.2354 -> 3101.4736
$ws.Cells.Item(116, 10).Value = 3024.818  # J116: 2918.5557 -> 3024.818
$ws.Cells.Item(116, 12).Value = 3024.818  # L116: 2918.5557 -> 3024.818
$ws.Cells.Item(116, 14).Value = -9908.817999999999  # N116: -9802.555700000001 -> -9908.817999999999
# Row 117
$ws.Cells.Item(117, 8).Value = 50000  # H117: 42600 -> 50000
$ws.Cells.Item(117, 10).Value = 50000  # J117: 42600 -> 50000
$ws.Cells.Item(117, 12).Value = 50000  # L117: 42600 -> 50000
$ws.Cells.Item(117, 14).Value = -59178  # N117: -51778 -> -59178
# Row 137
$ws.Cells.Item(137, 8).Value = 2306.64  # H137: 2376.6956 -> 2306.64
$ws.Cells.Item(137, 9).Value = 2283.35  # I137: 2377.1052 -> 2283.35
$ws.Cells.Item(137, 10).Value = 2399.8  # J137: 2374.75 -> 2399.8
$ws.Cells.Item(137, 11).Value = 6850.049999999999  # K137: 7131.3156 -> 6850.049999999999
$ws.Cells.Item(137, 12).Value = 7199.400000000001  # L137: 7124.25 -> 7199.400000000001
$ws.Cells.Item(137, 13).Value = -4300.049999999999  # M137: -4581.3156 -> -4300.049999999999
$ws.Cells.Item(137, 14).Value = -12299.4  # N137: -12224.25 -> -12299.4
# Row 138
$ws.Cells.Item(138, 8).Value = 6414745  # H138: 6540581 -> 6414745
$ws.Cells.Item(138, 9).Value = 1774.4  # I138: 2019.8948 -> 1774.4
$ws.Cells.Item(138, 11).Value = 5323.200000000001  # K138: 6059.6844 -> 5323.200000000001
$ws.Cells.Item(138, 13).Value = -183.2000000000007  # M138: -919.6844000000001 -> -183.2000000000007

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 14296847  # H32: 14717335 -> 14296847
$ws.Cells.Item(32, 9).Value = 23814668  # I32: 25005388 -> 23814668
$ws.Cells.Item(32, 11).Value = 23814668  # K32: 25005388 -> 23814668
$ws.Cells.Item(32, 13).Value = -23814381  # M32: -25005101 -> -23814381
# Row 61
$ws.Cells.Item(61, 8).Value = 41671090  # H61: 45458932 -> 41671090
$ws.Cells.Item(61, 9).Value = 55558356  # I61: 58826384 -> 55558356
$ws.Cells.Item(61, 10).Value = 9283.333000000001  # J61: 9600 -> 9283.333000000001
$ws.Cells.Item(61, 11).Value = 55558356  # K61: 58826384 -> 55558356
$ws.Cells.Item(61, 12).Value = 9283.333000000001  # L61: 9600 -> 9283.333000000001
$ws.Cells.Item(61, 13).Value = -55558144  # M61: -58826172 -> -55558144
$ws.Cells.Item(61, 14).Value = -9707.333000000001  # N61: -10024 -> -9707.333000000001
# Row 88
$ws.Cells.Item(88, 8).Value = 9017.134  # H88: 9591.857 -> 9017.134
$ws.Cells.Item(88, 9).Value = 12132.2  # I88: 13372.333 -> 12132.2
$ws.Cells.Item(88, 11).Value = 12132.2  # K88: 13372.333 -> 12132.2
$ws.Cells.Item(88, 13).Value = -11726.2  # M88: -12966.333 -> -11726.2
# Row 91
$ws.Cells.Item(91, 8).Value = 9017.134  # H91: 9591.857 -> 9017.134
$ws.Cells.Item(91, 9).Value = 12132.2  # I91: 13372.333 -> 12132.2
$ws.Cells.Item(91, 11).Value = 12132.2  # K91: 13372.333 -> 12132.2
$ws.Cells.Item(91, 13).Value = -10728.2  # M91: -11968.333 -> -10728.2
# Row 110
$ws.Cells.Item(110, 8).Value = 2223.3076  # H110: 2306.1785 -> 2223.3076
$ws.Cells.Item(110, 9).Value = 2173.8823  # I110: 2034.8948 -> 2173.8823
$ws.Cells.Item(110, 10).Value = 2316.6667  # J110: 2878.889 -> 2316.6667
$ws.Cells.Item(110, 11).Value = 2173.8823  # K110: 2034.8948 -> 2173.8823
$ws.Cells.Item(110, 12).Value = 2316.6667  # L110: 2878.889 -> 2316.6667
$ws.Cells.Item(110, 13).Value = -128.8823000000002  # M110: 10.10519999999997 -> -128.8823000000002
$ws.Cells.Item(110, 14).Value = -6406.6667  # N110: -6968.889 -> -6406.6667
# Row 132
$ws.Cells.Item(132, 8).Value = 33342042  # H132: 34491748 -> 33342042
$ws.Cells.Item(132, 9).Value = 10022.68  # I132: 10416.125 -> 10022.68
$ws.Cells.Item(132, 11).Value = 30068.04  # K132: 31248.375 -> 30068.04
$ws.Cells.Item(132, 13).Value = -27538.04  # M132: -28718.375 -> -27538.04
# Row 136
$ws.Cells.Item(136, 8).Value = 41671090  # H136: 45458932 -> 41671090
$ws.Cells.Item(136, 9).Value = 55558356  # I136: 58826384 -> 55558356
$ws.Cells.Item(136, 10).Value = 9283.333000000001  # J136: 9600 -> 9283.333000000001
$ws.Cells.Item(136, 11).Value = 166675068  # K136: 176479152 -> 166675068
$ws.Cells.Item(136, 12).Value = 27849.999  # L136: 28800 -> 27849.999
$ws.Cells.Item(136, 13).Value = -166672518  # M136: -176476602 -> -166672518
$ws.Cells.Item(136, 14).Value = -32949.999  # N136: -33900 -> -32949.999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 8
$ws.Cells.Item(8, 8).Value = 10687  # H8: 20147 -> 10687
$ws.Cells.Item(8, 9).Value = 10687  # I8: 26498.5 -> 10687
$ws.Cells.Item(8, 10).Value = 0  # J8: 7444 -> 0
$ws.Cells.Item(8, 11).Value = 10687  # K8: 26498.5 -> 10687
$ws.Cells.Item(8, 12).Value = 0  # L8: 7444 -> 0
$ws.Cells.Item(8, 13).ClearContents()  # M8: -26358.5 -> (removed)
$ws.Cells.Item(8, 14).Value = -10547  # N8: -7724 -> -10547
# Row 75
$ws.Cells.Item(75, 8).Value = 22982.363  # H75: 19682.6 -> 22982.363
$ws.Cells.Item(75, 9).Value = 5262.2856  # I75: 4976.5713 -> 5262.2856
$ws.Cells.Item(75, 10).Value = 53992.5  # J75: 53996.668 -> 53992.5
$ws.Cells.Item(75, 11).Value = 5262.2856  # K75: 4976.5713 -> 5262.2856
$ws.Cells.Item(75, 12).Value = 53992.5  # L75: 53996.668 -> 53992.5
$ws.Cells.Item(75, 13).Value = -4326.2856  # M75: -4040.5713 -> -4326.2856
$ws.Cells.Item(75, 14).Value = -55864.5  # N75: -55868.668 -> -55864.5
# Row 78
$ws.Cells.Item(78, 8).Value = 22982.363  # H78: 19682.6 -> 22982.363
$ws.Cells.Item(78, 9).Value = 5262.2856  # I78: 4976.5713 -> 5262.2856
$ws.Cells.Item(78, 10).Value = 53992.5  # J78: 53996.668 -> 53992.5
$ws.Cells.Item(78, 11).Value = 15786.8568  # K78: 14929.7139 -> 15786.8568
$ws.Cells.Item(78, 12).Value = 161977.5  # L78: 161990.004 -> 161977.5
$ws.Cells.Item(78, 13).Value = -11106.8568  # M78: -10249.7139 -> -11106.8568
$ws.Cells.Item(78, 14).Value = -171337.5  # N78: -171350.004 -> -171337.5
# Row 109
$ws.Cells.Item(109, 8).Value = 0  # H109: 69684 -> 0
$ws.Cells.Item(109, 10).Value = 0  # J109: 69684 -> 0
$ws.Cells.Item(109, 12).ClearContents()  # L109: 69684 -> (removed)
$ws.Cells.Item(109, 14).Value = 0  # N109: -72458 -> 0

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Cells.Item(2, 8).Value = 7949.5  # H2: 14999 -> 7949.5
$ws.Cells.Item(2, 9).Value = 7949.5  # I2: 14999 -> 7949.5
$ws.Cells.Item(2, 11).Value = 7949.5  # K2: 14999 -> 7949.5
$ws.Cells.Item(2, 13).Value = -7836.5  # M2: -14886 -> -7836.5
# Row 52
$ws.Cells.Item(52, 8).Value = 148972  # H52: 148977 -> 148972
$ws.Cells.Item(52, 10).Value = 148972  # J52: 148977 -> 148972
$ws.Cells.Item(52, 12).Value = 148972  # L52: 148977 -> 148972
$ws.Cells.Item(52, 14).Value = -149560  # N52: -149565 -> -149560

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Cells.Item(68, 8).Value = 1266.6666  # H68: 1113.2858 -> 1266.6666
$ws.Cells.Item(68, 9).Value = 1033.3334  # I68: 1096.6666 -> 1033.3334
$ws.Cells.Item(68, 10).Value = 1500  # J68: 1125.75 -> 1500
$ws.Cells.Item(68, 11).Value = 3100.0002  # K68: 3289.9998 -> 3100.0002
$ws.Cells.Item(68, 12).Value = 4500  # L68: 3377.25 -> 4500
$ws.Cells.Item(68, 13).Value = -2289.0002  # M68: -2478.9998 -> -2289.0002
$ws.Cells.Item(68, 14).Value = -6122  # N68: -4999.25 -> -6122
# Row 71
$ws.Cells.Item(71, 8).Value = 1266.6666  # H71: 1113.2858 -> 1266.6666
$ws.Cells.Item(71, 9).Value = 1033.3334  # I71: 1096.6666 -> 1033.3334
$ws.Cells.Item(71, 10).Value = 1500  # J71: 1125.75 -> 1500
$ws.Cells.Item(71, 11).Value = 9300.000599999999  # K71: 9869.999400000001 -> 9300.000599999999
$ws.Cells.Item(71, 12).Value = 13500  # L71: 10131.75 -> 13500
$ws.Cells.Item(71, 13).Value = -5244.000599999999  # M71: -5813.999400000001 -> -5244.000599999999
$ws.Cells.Item(71, 14).Value = -21612  # N71: -18243.75 -> -21612

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 4
$ws.Cells.Item(4, 8).Value = 1000  # H4: 0 -> 1000
$ws.Cells.Item(4, 9).Value = 1000  # I4: 0 -> 1000
$ws.Cells.Item(4, 11).Value = 1000  # K4: 0 -> 1000
$ws.Cells.Item(4, 13).Value = -888  # M4: None -> -888
# Row 24
$ws.Cells.Item(24, 8).Value = 23571.857  # H24: 30006 -> 23571.857
$ws.Cells.Item(24, 9).Value = 20001.5  # I24: 30006 -> 20001.5
$ws.Cells.Item(24, 10).Value = 25000  # J24: 0 -> 25000
$ws.Cells.Item(24, 11).Value = 20001.5  # K24: 30006 -> 20001.5
$ws.Cells.Item(24, 12).Value = 25000  # L24: 0 -> 25000
$ws.Cells.Item(24, 13).Value = -19828.5  # M24: -29833 -> -19828.5
$ws.Cells.Item(24, 14).Value = -25346  # N24: None -> -25346
# Row 80
$ws.Cells.Item(80, 8).Value = 3358.8572  # H80: 3357.4375 -> 3358.8572
$ws.Cells.Item(80, 10).Value = 3822.8333  # J80: 3704 -> 3822.8333
$ws.Cells.Item(80, 12).Value = 3822.8333  # L80: 3704 -> 3822.8333
$ws.Cells.Item(80, 14).Value = -5818.8333  # N80: -5700 -> -5818.8333
# Row 83
$ws.Cells.Item(83, 8).Value = 3358.8572  # H83: 3357.4375 -> 3358.8572
$ws.Cells.Item(83, 10).Value = 3822.8333  # J83: 3704 -> 3822.8333
$ws.Cells.Item(83, 12).Value = 19114.1665  # L83: 18520 -> 19114.1665
$ws.Cells.Item(83, 14).Value = -29098.1665  # N83: -28504 -> -29098.1665

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Cells.Item(2, 8).Value = 10  # H2: 4 -> 10
$ws.Cells.Item(2, 9).Value = 10  # I2: 4 -> 10
$ws.Cells.Item(2, 11).Value = 10  # K2: 4 -> 10
$ws.Cells.Item(2, 13).Value = 102  # M2: 108 -> 102
# Row 7
$ws.Cells.Item(7, 8).Value = 4130.6665  # H7: 4124.5684 -> 4130.6665
$ws.Cells.Item(7, 10).Value = 4564.0835  # J7: 4571.2607 -> 4564.0835
$ws.Cells.Item(7, 12).Value = 4564.0835  # L7: 4571.2607 -> 4564.0835
$ws.Cells.Item(7, 14).Value = -4788.0835  # N7: -4795.2607 -> -4788.0835
# Row 55
$ws.Cells.Item(55, 8).Value = 607.1905  # H55: 580.5 -> 607.1905
$ws.Cells.Item(55, 9).Value = 387.25  # I55: 365.64706 -> 387.25
$ws.Cells.Item(55, 11).Value = 387.25  # K55: 365.64706 -> 387.25
$ws.Cells.Item(55, 13).Value = -214.25  # M55: -192.64706 -> -214.25
# Row 61
$ws.Cells.Item(61, 8).Value = 2158.0667  # H61: 1903.0952 -> 2158.0667
$ws.Cells.Item(61, 9).Value = 831  # I61: 1044.5454 -> 831
$ws.Cells.Item(61, 10).Value = 5807.5  # J61: 2847.5 -> 5807.5
$ws.Cells.Item(61, 11).Value = 831  # K61: 1044.5454 -> 831
$ws.Cells.Item(61, 12).Value = 5807.5  # L61: 2847.5 -> 5807.5
$ws.Cells.Item(61, 13).Value = -629  # M61: -842.5454 -> -629
$ws.Cells.Item(61, 14).Value = -6211.5  # N61: -3251.5 -> -6211.5
# Row 100
$ws.Cells.Item(100, 8).Value = 2590.6365  # H100: 2780.2 -> 2590.6365
$ws.Cells.Item(100, 9).Value = 1682.6666  # I100: 1880.2 -> 1682.6666
$ws.Cells.Item(100, 11).Value = 1682.6666  # K100: 1880.2 -> 1682.6666
$ws.Cells.Item(100, 13).Value = -1141.6666  # M100: -1339.2 -> -1141.6666
# Row 113
$ws.Cells.Item(113, 8).Value = 2158.0667  # H113: 1903.0952 -> 2158.0667
$ws.Cells.Item(113, 9).Value = 831  # I113: 1044.5454 -> 831
$ws.Cells.Item(113, 10).Value = 5807.5  # J113: 2847.5 -> 5807.5
$ws.Cells.Item(113, 11).Value = 831  # K113: 1044.5454 -> 831
$ws.Cells.Item(113, 12).Value = 5807.5  # L113: 2847.5 -> 5807.5
$ws.Cells.Item(113, 13).Value = 1339  # M113: 1125.4546 -> 1339
$ws.Cells.Item(113, 14).Value = -10147.5  # N113: -7187.5 -> -10147.5
# Row 126
$ws.Cells.Item(126, 8).Value = 4130.6665  # H126: 4124.5684 -> 4130.6665
$ws.Cells.Item(126, 10).Value = 4564.0835  # J126: 4571.2607 -> 4564.0835
$ws.Cells.Item(126, 12).Value = 13692.2505  # L126: 13713.7821 -> 13692.2505
$ws.Cells.Item(126, 14).Value = -18632.2505  # N126: -18653.7821 -> -18632.2505
# Row 132
$ws.Cells.Item(132, 8).Value = 71430390  # H132: 64517892 -> 71430390
$ws.Cells.Item(132, 9).Value = 1783.5454  # I132: 1758.1538 -> 1783.5454
$ws.Cells.Item(132, 10).Value = 333335300  # J132: 400001800 -> 333335300
$ws.Cells.Item(132, 11).Value = 5350.6362  # K132: 5274.4614 -> 5350.6362
$ws.Cells.Item(132, 12).Value = 1000005900  # L132: 1200005400 -> 1000005900
$ws.Cells.Item(132, 13).Value = -2820.6362  # M132: -2744.4614 -> -2820.6362
$ws.Cells.Item(132, 14).Value = -1000010960  # N132: -1200010460 -> -1000010960
# Row 136
$ws.Cells.Item(136, 8).Value = 2526.9333  # H136: 2562.0364 -> 2526.9333
$ws.Cells.Item(136, 9).Value = 1928.6227  # I136: 1906.5209 -> 1928.6227
$ws.Cells.Item(136, 11).Value = 5785.8681  # K136: 5719.5627 -> 5785.8681
$ws.Cells.Item(136, 13).Value = -3235.8681  # M136: -3169.5627 -> -3235.8681

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 20
$ws.Cells.Item(20, 8).Value = 0  # H20: 3000 -> 0
$ws.Cells.Item(20, 9).Value = 0  # I20: 3000 -> 0
$ws.Cells.Item(20, 11).Value = 0  # K20: 3000 -> 0
$ws.Cells.Item(20, 13).ClearContents()  # M20: -2760 -> (removed)
# Row 81
$ws.Cells.Item(81, 8).Value = 1000  # H81: 800 -> 1000
$ws.Cells.Item(81, 9).Value = 750  # I81: 500 -> 750
$ws.Cells.Item(81, 10).Value = 2000  # J81: 1400 -> 2000
$ws.Cells.Item(81, 11).Value = 1500  # K81: 1000 -> 1500
$ws.Cells.Item(81, 12).Value = 4000  # L81: 2800 -> 4000
$ws.Cells.Item(81, 13).Value = -439  # M81: 61 -> -439
$ws.Cells.Item(81, 14).Value = -6122  # N81: -4922 -> -6122
# Row 84
$ws.Cells.Item(84, 8).Value = 1000  # H84: 800 -> 1000
$ws.Cells.Item(84, 9).Value = 750  # I84: 500 -> 750
$ws.Cells.Item(84, 10).Value = 2000  # J84: 1400 -> 2000
$ws.Cells.Item(84, 11).Value = 7500  # K84: 5000 -> 7500
$ws.Cells.Item(84, 12).Value = 20000  # L84: 14000 -> 20000
$ws.Cells.Item(84, 13).Value = -2196  # M84: 304 -> -2196
$ws.Cells.Item(84, 14).Value = -30608  # N84: -24608 -> -30608
# Row 113
$ws.Cells.Item(113, 8).Value = 459.4  # H113: 457.75 -> 459.4
$ws.Cells.Item(113, 10).Value = 1140.1666  # J113: 1268.4 -> 1140.1666
$ws.Cells.Item(113, 12).Value = 3420.4998  # L113: 3805.2 -> 3420.4998
$ws.Cells.Item(113, 14).Value = -7760.4998  # N113: -8145.200000000001 -> -7760.4998
# Row 116
$ws.Cells.Item(116, 8).Value = 76000  # H116: 44500 -> 76000
$ws.Cells.Item(116, 10).Value = 76000  # J116: 44500 -> 76000
$ws.Cells.Item(116, 12).Value = 76000  # L116: 44500 -> 76000
$ws.Cells.Item(116, 14).Value = -85178  # N116: -53678 -> -85178
# Row 136
$ws.Cells.Item(136, 8).Value = 1670.1562  # H136: 1711.1613 -> 1670.1562
$ws.Cells.Item(136, 9).Value = 1648.9  # I136: 1692 -> 1648.9
$ws.Cells.Item(136, 11).Value = 4946.700000000001  # K136: 5076 -> 4946.700000000001
$ws.Cells.Item(136, 13).Value = -2396.700000000001  # M136: -2526 -> -2396.700000000001
